# Natmi following Dr Hou advice
# Update the Ncam1-Gfra1 LR-pairs sheet: sending-cluster "sCs" becomes "M2",
# and a new set of rows is added for the (now re-labelled) "sCs" sending
# cluster, with refreshed statistics throughout (3 replicate samples instead
# of 1, new weighted expression values, etc.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, "A").Value = "ECs"
$ws.Cells.Item(2, "B").Value = "Ncam1"
$ws.Cells.Item(2, "C").Value = "Gfra1"
$ws.Cells.Item(2, "D").Value = "FAPs"
$ws.Cells.Item(2, "E").Value = 3
$ws.Cells.Item(2, "F").Value = 1
$ws.Cells.Item(2, "G").Value = 0.300794
$ws.Cells.Item(2, "H").Value = 0.902382
$ws.Cells.Item(2, "I").Value = 0.003836667885433928
$ws.Cells.Item(2, "J").Value = 0.003836667885433928
$ws.Cells.Item(2, "K").Value = 3
$ws.Cells.Item(2, "L").Value = 1
$ws.Cells.Item(2, "M").Value = 15.933008
$ws.Cells.Item(2, "N").Value = 47.799024
$ws.Cells.Item(2, "O").Value = 0.7524552897049799
$ws.Cells.Item(2, "P").Value = 0.7524552897049799
$ws.Cells.Item(2, "Q").Value = 4.792553208352
$ws.Cells.Item(2, "R").Value = 43.132978875168
$ws.Cells.Item(2, "S").Value = 0.002886921045235979
$ws.Cells.Item(2, "T").Value = 0.002886921045235979

# Row 3
$ws.Cells.Item(3, "A").Value = "ECs"
$ws.Cells.Item(3, "B").Value = "Ncam1"
$ws.Cells.Item(3, "C").Value = "Gfra1"
$ws.Cells.Item(3, "D").Value = "sCs"
$ws.Cells.Item(3, "E").Value = 3
$ws.Cells.Item(3, "F").Value = 1
$ws.Cells.Item(3, "G").Value = 0.300794
$ws.Cells.Item(3, "H").Value = 0.902382
$ws.Cells.Item(3, "I").Value = 0.003836667885433928
$ws.Cells.Item(3, "J").Value = 0.003836667885433928
$ws.Cells.Item(3, "K").Value = 3
$ws.Cells.Item(3, "L").Value = 1
$ws.Cells.Item(3, "M").Value = 5.241682666666667
$ws.Cells.Item(3, "N").Value = 15.725048
$ws.Cells.Item(3, "O").Value = 0.2475447102950201
$ws.Cells.Item(3, "P").Value = 0.2475447102950201
$ws.Cells.Item(3, "Q").Value = 1.576666696037333
$ws.Cells.Item(3, "R").Value = 14.190000264336
$ws.Cells.Item(3, "S").Value = 0.0009497468401979491
$ws.Cells.Item(3, "T").Value = 0.0009497468401979491

# Row 4
$ws.Cells.Item(4, "A").Value = "FAPs"
$ws.Cells.Item(4, "B").Value = "Ncam1"
$ws.Cells.Item(4, "C").Value = "Gfra1"
$ws.Cells.Item(4, "D").Value = "FAPs"
$ws.Cells.Item(4, "E").Value = 3
$ws.Cells.Item(4, "F").Value = 1
$ws.Cells.Item(4, "G").Value = 6.268658666666667
$ws.Cells.Item(4, "H").Value = 18.805976
$ws.Cells.Item(4, "I").Value = 0.07995758356598558
$ws.Cells.Item(4, "J").Value = 0.07995758356598558
$ws.Cells.Item(4, "K").Value = 3
$ws.Cells.Item(4, "L").Value = 1
$ws.Cells.Item(4, "M").Value = 15.933008
$ws.Cells.Item(4, "N").Value = 47.799024
$ws.Cells.Item(4, "O").Value = 0.7524552897049799
$ws.Cells.Item(4, "P").Value = 0.7524552897049799
$ws.Cells.Item(4, "Q").Value = 99.87858868526935
$ws.Cells.Item(4, "R").Value = 898.9072981674241
$ws.Cells.Item(4, "S").Value = 0.06016450670625381
$ws.Cells.Item(4, "T").Value = 0.06016450670625381

# Row 5
$ws.Cells.Item(5, "A").Value = "FAPs"
$ws.Cells.Item(5, "B").Value = "Ncam1"
$ws.Cells.Item(5, "C").Value = "Gfra1"
$ws.Cells.Item(5, "D").Value = "sCs"
$ws.Cells.Item(5, "E").Value = 3
$ws.Cells.Item(5, "F").Value = 1
$ws.Cells.Item(5, "G").Value = 6.268658666666667
$ws.Cells.Item(5, "H").Value = 18.805976
$ws.Cells.Item(5, "I").Value = 0.07995758356598558
$ws.Cells.Item(5, "J").Value = 0.07995758356598558
$ws.Cells.Item(5, "K").Value = 3
$ws.Cells.Item(5, "L").Value = 1
$ws.Cells.Item(5, "M").Value = 5.241682666666667
$ws.Cells.Item(5, "N").Value = 15.725048
$ws.Cells.Item(5, "O").Value = 0.2475447102950201
$ws.Cells.Item(5, "P").Value = 0.2475447102950201
$ws.Cells.Item(5, "Q").Value = 32.85831947631645
$ws.Cells.Item(5, "R").Value = 295.724875286848
$ws.Cells.Item(5, "S").Value = 0.01979307685973176
$ws.Cells.Item(5, "T").Value = 0.01979307685973176

# Row 6
$ws.Cells.Item(6, "A").Value = "M2"
$ws.Cells.Item(6, "B").Value = "Ncam1"
$ws.Cells.Item(6, "C").Value = "Gfra1"
$ws.Cells.Item(6, "D").Value = "FAPs"
$ws.Cells.Item(6, "E").Value = 1
$ws.Cells.Item(6, "F").Value = 0.3333333333333333
$ws.Cells.Item(6, "G").Value = 0.04541633333333334
$ws.Cells.Item(6, "H").Value = 0.136249
$ws.Cells.Item(6, "I").Value = 0.0005792914339187697
$ws.Cells.Item(6, "J").Value = 0.0005792914339187697
$ws.Cells.Item(6, "K").Value = 3
$ws.Cells.Item(6, "L").Value = 1
$ws.Cells.Item(6, "M").Value = 15.933008
$ws.Cells.Item(6, "N").Value = 47.799024
$ws.Cells.Item(6, "O").Value = 0.7524552897049799
$ws.Cells.Item(6, "P").Value = 0.7524552897049799
$ws.Cells.Item(6, "Q").Value = 0.7236188023306668
$ws.Cells.Item(6, "R").Value = 6.512569220976001
$ws.Cells.Item(6, "S").Value = 0.0004358909037329611
$ws.Cells.Item(6, "T").Value = 0.0004358909037329611

# Row 7
$ws.Cells.Item(7, "A").Value = "M2"
$ws.Cells.Item(7, "B").Value = "Ncam1"
$ws.Cells.Item(7, "C").Value = "Gfra1"
$ws.Cells.Item(7, "D").Value = "sCs"
$ws.Cells.Item(7, "E").Value = 1
$ws.Cells.Item(7, "F").Value = 0.3333333333333333
$ws.Cells.Item(7, "G").Value = 0.04541633333333334
$ws.Cells.Item(7, "H").Value = 0.136249
$ws.Cells.Item(7, "I").Value = 0.0005792914339187697
$ws.Cells.Item(7, "J").Value = 0.0005792914339187697
$ws.Cells.Item(7, "K").Value = 3
$ws.Cells.Item(7, "L").Value = 1
$ws.Cells.Item(7, "M").Value = 5.241682666666667
$ws.Cells.Item(7, "N").Value = 15.725048
$ws.Cells.Item(7, "O").Value = 0.2475447102950201
$ws.Cells.Item(7, "P").Value = 0.2475447102950201
$ws.Cells.Item(7, "Q").Value = 0.2380580072168889
$ws.Cells.Item(7, "R").Value = 2.142522064952
$ws.Cells.Item(7, "S").Value = 0.0001434005301858087
$ws.Cells.Item(7, "T").Value = 0.0001434005301858087

# Row 8
$ws.Cells.Item(8, "A").Value = "sCs"
$ws.Cells.Item(8, "B").Value = "Ncam1"
$ws.Cells.Item(8, "C").Value = "Gfra1"
$ws.Cells.Item(8, "D").Value = "FAPs"
$ws.Cells.Item(8, "E").Value = 3
$ws.Cells.Item(8, "F").Value = 1
$ws.Cells.Item(8, "G").Value = 71.78493233333333
$ws.Cells.Item(8, "H").Value = 215.354797
$ws.Cells.Item(8, "I").Value = 0.9156264571146617
$ws.Cells.Item(8, "J").Value = 0.9156264571146617
$ws.Cells.Item(8, "K").Value = 3
$ws.Cells.Item(8, "L").Value = 1
$ws.Cells.Item(8, "M").Value = 15.933008
$ws.Cells.Item(8, "N").Value = 47.799024
$ws.Cells.Item(8, "O").Value = 0.7524552897049799
$ws.Cells.Item(8, "P").Value = 0.7524552897049799
$ws.Cells.Item(8, "Q").Value = 1143.749901146459
$ws.Cells.Item(8, "R").Value = 10293.74911031813
$ws.Cells.Item(8, "S").Value = 0.6889679710497572
$ws.Cells.Item(8, "T").Value = 0.6889679710497572

# Row 9
$ws.Cells.Item(9, "A").Value = "sCs"
$ws.Cells.Item(9, "B").Value = "Ncam1"
$ws.Cells.Item(9, "C").Value = "Gfra1"
$ws.Cells.Item(9, "D").Value = "sCs"
$ws.Cells.Item(9, "E").Value = 3
$ws.Cells.Item(9, "F").Value = 1
$ws.Cells.Item(9, "G").Value = 71.78493233333333
$ws.Cells.Item(9, "H").Value = 215.354797
$ws.Cells.Item(9, "I").Value = 0.9156264571146617
$ws.Cells.Item(9, "J").Value = 0.9156264571146617
$ws.Cells.Item(9, "K").Value = 3
$ws.Cells.Item(9, "L").Value = 1
$ws.Cells.Item(9, "M").Value = 5.241682666666667
$ws.Cells.Item(9, "N").Value = 15.725048
$ws.Cells.Item(9, "O").Value = 0.2475447102950201
$ws.Cells.Item(9, "P").Value = 0.2475447102950201
$ws.Cells.Item(9, "Q").Value = 376.2738355394729
$ws.Cells.Item(9, "R").Value = 3386.464519855256
$ws.Cells.Item(9, "S").Value = 0.2266584860649046
$ws.Cells.Item(9, "T").Value = 0.2266584860649046
